$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (C11): part changed from Samsung CL21F104ZBCNNNC (Y5V dielectric) to
# Yageo CC0805KRX7R9BB104 (X7R dielectric).
$ws.Range("D3").Value = "CAP CER 0.1UF 50V X7R 0805"
$ws.Range("E3").Value = "X7R"
$ws.Range("G3").Value = "Yageo"
$ws.Range("I3").Value = "CC0805KRX7R9BB104"
$ws.Range("J3").Value = "311-1140-1-ND"
$ws.Range("K3").Value = "603-CC805KRX7R9BB104"
$ws.Range("L3").Value = 0.1

# Row 29 (IC7 MC33814): drop stray trailing newline baked into the Mouser P/N string.
$ws.Range("K29").Value = "841-MC33814AE"
